$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-6: 45184 -> 45185
$ws.Range("C2:C6").Value = 45185
